$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 276
$ws1.Range("F4").Value = 1052
$ws1.Range("F5").Value = 559

# Sheet "全部类型" (All Types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 276
$ws4.Range("F4").Value = 1052
$ws4.Range("F6").Value = 559
